$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new data rows right after the header row (row 1), before the
# existing row 11 ("Especial" for date 45145). This pushes the old
# rows 11..94 down to 15..98, matching the target diff.
$ws.Rows("11:14").Insert()

# Fixed values shared by every data row in this sheet.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = '$/bandeja 3 kilos'
$origen = "Región de Arica y Parinacota"
$kgUnidad = 3

# New rows for the latest reporting date (45149 => 2023-08-11), one per
# quality grade, matching the commit's weekly price update.
$newRows = @(
    @{ Row = 11; Calidad = "Especial"; Volumen = 60;  Min = 7000; Max = 8000; Prom = 7500; Kg = 2500 },
    @{ Row = 12; Calidad = "Primera";  Volumen = 80;  Min = 5000; Max = 6000; Prom = 5500; Kg = 1833 },
    @{ Row = 13; Calidad = "Segunda";  Volumen = 100; Min = 4000; Max = 5000; Prom = 4500; Kg = 1500 },
    @{ Row = 14; Calidad = "Tercera";  Volumen = 160; Min = 2000; Max = 3000; Prom = 2500; Kg = 833 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $mercadoId
    $ws.Cells.Item($row, 2).Value2 = $mercado
    $ws.Cells.Item($row, 3).Value2 = $region
    $ws.Cells.Item($row, 4).Value2 = 45149
    $ws.Cells.Item($row, 5).Value2 = $codreg
    $ws.Cells.Item($row, 6).Value2 = $tipo
    $ws.Cells.Item($row, 7).Value2 = $productoId
    $ws.Cells.Item($row, 8).Value2 = $producto
    $ws.Cells.Item($row, 9).Value2 = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.Min
    $ws.Cells.Item($row, 15).Value2 = $r.Max
    $ws.Cells.Item($row, 16).Value2 = $r.Prom
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.Kg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
